# Auto-generated edit script: update crypto price/volume columns (D, E)
# to reflect the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.011.30"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "3.527.08"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "600.23"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").Value = "183.55"
$ws.Range("E6").Value = "  +5.66%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("D10").Value = "7.14"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("E11").Value = "  +1.35%  "
$ws.Range("D12").Value = "4.141.43"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "32.49"
$ws.Range("E13").Value = "  +12.15%  "
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "67.995.52"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "0.0000182"
$ws.Range("E16").Value = "  +0.79%  "
$ws.Range("D17").Value = "3.536.26"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "6.41"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").Value = "14.96"
$ws.Range("E19").Value = "  +4.83%  "
$ws.Range("D20").Value = "399.34"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "8.13"
$ws.Range("E21").Value = "  +1.78%  "
$ws.Range("D22").Value = "73.77"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("D23").Value = "0.546"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "0.0000125"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D26").NumberFormat = "@"  # preserve trailing zero as text
$ws.Range("D26").Value = "5.70"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "10.77"
$ws.Range("E27").Value = "  +5.48%  "
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").Value = "6.31"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("E32").Value = "  +0.96%  "
$ws.Range("D33").Value = "24.11"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").NumberFormat = "@"  # preserve trailing zero as text
$ws.Range("D34").Value = "7.50"
$ws.Range("E34").Value = "  +1.30%  "
$ws.Range("D36").Value = "1.68"
$ws.Range("E36").Value = "  +2.75%  "
$ws.Range("D37").Value = "163.65"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "1.97"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("D39").NumberFormat = "@"  # preserve trailing zero as text
$ws.Range("D39").Value = "0.880"
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("D40").Value = "7.16"
$ws.Range("E40").Value = "  +3.70%  "
$ws.Range("E41").Value = "  +7.21%  "
$ws.Range("D42").Value = "4.77"
$ws.Range("E42").Value = "  +1.84%  "
$ws.Range("D43").Value = "27.17"
$ws.Range("E43").Value = "  +2.69%  "
$ws.Range("D44").Value = "27.63"
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").Value = "2.881.38"
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("D46").Value = "0.0742"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("E47").Value = "  -1.03%  "
$ws.Range("D48").Value = "351.49"
$ws.Range("E48").Value = "  +3.62%  "
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").NumberFormat = "@"  # preserve trailing zero as text
$ws.Range("D51").Value = "33.90"
$ws.Range("E51").Value = "  +1.19%  "
